$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 31: update B31/C31 from "OP11"/"LogIn" to the new SSO WAYFLess test entry ---
$ws.Range("B31").Value = "OPQA-5682"
$ws.Range("C31").Value = "Verify the SSO authentication via WAYFLess/direct URL"

# --- Insert 6 new rows (32:37) for the new DRAIAM112-DRAIAM117 test cases ---
$ws.Rows("32:37").Insert()

# Copy the formatting of row 31 (border/wrap, no fill) onto the freshly inserted rows
$ws.Range("A31:E31").Copy()
$ws.Range("A32:E37").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A32").Value = "DRAIAM112"
$ws.Range("B32").Value = "OPQA-5689"
$ws.Range("C32").Value = 'verify that "Your account has been suspended Questions? Contact [Application_customer_care]" error message when user is an existing STeAM User that is suspended.'
$ws.Range("D32").Value = "Y"

$ws.Range("A33").Value = "DRAIAM113"
$ws.Range("B33").Value = "OPQA-5694"
$ws.Range("C33").Value = "Verify that user able to login DRA ,If user exist in STeAM but not associated any DRA SSO claimticket"
$ws.Range("D33").Value = "Y"

$ws.Range("A34").Value = "DRAIAM114"
$ws.Range("B34").Value = "OPQA-5686"
$ws.Range("C34").Value = "Verify that user able to login DRA ,If user doesn't exist in STeAM but exist in IdP."
$ws.Range("D34").Value = "Y"

$ws.Range("A35").Value = "DRAIAM115"
$ws.Range("B35").Value = "OPQA-5683"
$ws.Range("C35").Value = "Verify that user with both a TD and DD subscription shall be taken to the interoperability page after successful authentication"
$ws.Range("D35").Value = "Y"

$ws.Range("A36").Value = "DRAIAM116"
$ws.Range("B36").Value = "OPQA-5684"
$ws.Range("C36").Value = "Verify that user with a TD-only subscription shall be taken to the TD homepage after successful authentication"
$ws.Range("D36").Value = "Y"

$ws.Range("A37").Value = "DRAIAM117"
$ws.Range("B37").Value = "OPQA-5685"
$ws.Range("C37").Value = "Verify that user with a DD-only subscription shall be taken to the DD homepage after successful authentication"
$ws.Range("D37").Value = "Y"

# --- The former blank rows (old 40:41) shifted down to 46:47 by the insert; remove them ---
$ws.Rows("46:47").Delete()

# --- Match the new selection shown in the target workbook ---
$ws.Range("A31").Select() | Out-Null
